$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$todos = @(
    "hidden objects should not be included in shade analysis?",
    "remove checkboxes to show/hide objects",
    "zoom to fit, 'F' key",
    "clear properties on de-select object",
    "object list sorted by group",
    "active area polygon rotation",
    "composite objects: regular fixed array"
)

$row = 67
foreach ($todo in $todos) {
    $ws.Cells.Item($row, 1).Value = "Not done"
    $ws.Cells.Item($row, 2).Value = $todo
    $row++
}

$ws.Range("B67:B73").Interior.ThemeColor = 9
$ws.Range("B67:B73").Interior.TintAndShade = 0.6

$ws.Range("F79").Select()
